# aggiornamento 15, 16, 17 marzo
# Append three new rows (227-229) of data below the existing table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: date-serial (col A), nuovi pos. (col B), somma mobile 7gg. (col C),
# somma mobile 7gg. per 100mila abitanti (col D)
$newRows = @(
    @{ Row = 227; A = 44301; B = 6; C = 24; D = 212.5963327132607 },
    @{ Row = 228; A = 44302; B = 6; C = 30; D = 265.7454158915759 },
    @{ Row = 229; A = 44303; B = 4; C = 28; D = 248.0290548321375 }
)

foreach ($r in $newRows) {
    $rowIdx = $r.Row
    $prevRowIdx = $rowIdx - 1

    # Copy formatting (style) of column A from the row above, which carries
    # the date number-format style (s="2") used throughout the column.
    $ws.Range("A$prevRowIdx").Copy($ws.Range("A$rowIdx"))

    $ws.Cells.Item($rowIdx, 1).Value = $r.A
    $ws.Cells.Item($rowIdx, 2).Value = $r.B
    $ws.Cells.Item($rowIdx, 3).Value = $r.C
    $ws.Cells.Item($rowIdx, 4).Value = $r.D
}

Write-Output "Added rows 227-229"
